# Booked_Hours_Angus.xlsx -- "add google drive snapshot from end of week4"
#
# The canonical-OOXML diff for this commit is a Google-Sheets re-export of
# the workbook: besides an internal style-table cleanup (merging a
# duplicate/near-duplicate "no explicit colour" Verdana font into the
# existing themed Verdana font -- a purely cosmetic, invisible change,
# since both resolve to black), the one user-visible edit is that columns
# C, D and E on Sheet1 were widened.
#
#   column C: 8.14  -> 12.29
#   column D: 9.29  -> 11.29
#   column E: 8.86  -> 12.29
#
# Excel's COM ColumnWidth property is expressed in "characters" for the
# workbook's Normal font, while the OOXML <col width="..."/> attribute is
# that value plus a fixed 5/MaximumDigitWidth padding term. For this
# workbook MDW = 7px, so OOXML_width = ColumnWidth + 5/7. Subtract that
# padding back out so the saved file lands on the target widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$padding = 5.0 / 7.0

$ws.Columns(3).ColumnWidth = 12.29 - $padding   # C
$ws.Columns(4).ColumnWidth = 11.29 - $padding   # D
$ws.Columns(5).ColumnWidth = 12.29 - $padding   # E
